$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds weekly price records for "Vega Modelo de Temuco - Achicoria".
# A new week's record was added (row 84 becomes the newest), which pushes the
# previously-newest three rows (84-86) down by one row, and what used to be
# row 86 now lands in a brand-new row 87.
#
# Before -> After:
#   Row84: D=45015 J=150            -> D=45041 J=100   (new record)
#   Row85: D=44727 J=35             -> D=45015 J=150
#   Row86: D=45007 J=25 O=Maule     -> D=44727 J=35  O=Región Metropolitana
#   Row87: (new row)                -> D=45007 J=25  O=Región del Maule (rest identical to old row86)

# Create the new row 87, carrying what used to be row 86's full content.
$ws.Range("A87").Value = 10
$ws.Range("B87").Value = "Vega Modelo de Temuco"
$ws.Range("C87").Value = "La Araucanía"
$ws.Range("D87").Value = 45007
$ws.Range("D87").NumberFormat = $ws.Range("D86").NumberFormat
$ws.Range("E87").Value = 9
$ws.Range("F87").Value = 100112010
$ws.Range("G87").Value = "Achicoria"
$ws.Range("H87").Value = "Sin especificar"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 25
$ws.Range("K87").Value = 10000
$ws.Range("L87").Value = 10000
$ws.Range("M87").Value = 10000
$ws.Range("N87").Value = "$/caja 18 unidades"
$ws.Range("O87").Value = "Región del Maule"
$ws.Range("P87").Value = 556
$ws.Range("Q87").Value = 18
$ws.Range("R87").Value = "Hortaliza"

# Shift the Date/Volumen values in rows 84-86 down one slot, and update
# row 86's Origen to its new value.
$ws.Range("D86").Value = 44727
$ws.Range("J86").Value = 35
$ws.Range("O86").Value = "Región Metropolitana"

$ws.Range("D85").Value = 45015
$ws.Range("J85").Value = 150

$ws.Range("D84").Value = 45041
$ws.Range("J84").Value = 100
